$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(657, 1).Value = "'2024-09-02"
$ws.Cells.Item(657, 3).Value = 2232.75
$ws.Cells.Item(657, 4).Value = 1537.550048828125
$ws.Cells.Item(657, 5).Value = 1687.900024414062
$ws.Cells.Item(657, 6).Value = 1970.599975585938
$ws.Cells.Item(657, 7).Value = 1111.550048828125
$ws.Cells.Item(657, 8).Value = 8540.35009765625
$ws.Cells.Item(657, 9).Value = 0
$ws.Cells.Item(657, 10).Value = 174.8915051884381

$ws.Cells.Item(658, 1).Value = "'2024-09-03"
$ws.Cells.Item(658, 3).Value = 2240.25
$ws.Cells.Item(658, 4).Value = 1530.599975585938
$ws.Cells.Item(658, 5).Value = 1687.5
$ws.Cells.Item(658, 6).Value = 1924.650024414062
$ws.Cells.Item(658, 7).Value = 1114
$ws.Cells.Item(658, 8).Value = 8497
$ws.Cells.Item(658, 9).Value = -0.005075915759957742
$ws.Cells.Item(658, 10).Value = 174.0037706409694

$ws.Cells.Item(659, 1).Value = "'2024-09-04"
$ws.Cells.Item(659, 3).Value = 2277.25
$ws.Cells.Item(659, 4).Value = 1556.550048828125
$ws.Cells.Item(659, 5).Value = 1686.550048828125
$ws.Cells.Item(659, 6).Value = 1924.650024414062
$ws.Cells.Item(659, 7).Value = 1127.900024414062
$ws.Cells.Item(659, 8).Value = 8572.900146484375
$ws.Cells.Item(659, 9).Value = 0.008932581674046723
$ws.Cells.Item(659, 10).Value = 175.5580735338119

$ws.Cells.Item(660, 1).Value = "'2024-09-05"
$ws.Cells.Item(660, 3).Value = 2290.199951171875
$ws.Cells.Item(660, 4).Value = 1555.75
$ws.Cells.Item(660, 5).Value = 1709.449951171875
$ws.Cells.Item(660, 6).Value = 1933.599975585938
$ws.Cells.Item(660, 7).Value = 1115.150024414062
$ws.Cells.Item(660, 8).Value = 8604.14990234375
$ws.Cells.Item(660, 9).Value = 0.00364517903223101
$ws.Cells.Item(660, 10).Value = 176.1980141423962

$ws.Cells.Item(661, 1).Value = "'2024-09-06"
$ws.Cells.Item(661, 3).Value = 2256.5
$ws.Cells.Item(661, 4).Value = 1559.900024414062
$ws.Cells.Item(661, 5).Value = 1702.699951171875
$ws.Cells.Item(661, 6).Value = 1928.400024414062
$ws.Cells.Item(661, 7).Value = 1100
$ws.Cells.Item(661, 8).Value = 8547.5
$ws.Cells.Item(661, 9).Value = -0.006584020848860234
$ws.Cells.Item(661, 10).Value = 175.0379227437549

$ws.Cells.Item(662, 1).Value = "'2024-09-09"
$ws.Cells.Item(662, 3).Value = 2216.800048828125
$ws.Cells.Item(662, 4).Value = 1546.25
$ws.Cells.Item(662, 5).Value = 1704.199951171875
$ws.Cells.Item(662, 6).Value = 1937.099975585938
$ws.Cells.Item(662, 7).Value = 1104.150024414062
$ws.Cells.Item(662, 8).Value = 8508.5
$ws.Cells.Item(662, 9).Value = -0.004562737642585551
$ws.Cells.Item(662, 10).Value = 174.239270624772

$ws.Cells.Item(663, 1).Value = "'2024-09-10"
$ws.Cells.Item(663, 3).Value = 2222.550048828125
$ws.Cells.Item(663, 4).Value = 1545.550048828125
$ws.Cells.Item(663, 5).Value = 1727.849975585938
$ws.Cells.Item(663, 6).Value = 1912.150024414062
$ws.Cells.Item(663, 7).Value = 1113.199951171875
$ws.Cells.Item(663, 8).Value = 8521.300048828125
$ws.Cells.Item(663, 9).Value = 0.001504383713712758
$ws.Cells.Item(663, 10).Value = 174.5013933457891

$ws.Cells.Item(664, 1).Value = "'2024-09-11"
$ws.Cells.Item(664, 3).Value = 2209.39990234375
$ws.Cells.Item(664, 4).Value = 1591.949951171875
$ws.Cells.Item(664, 5).Value = 1725.650024414062
$ws.Cells.Item(664, 6).Value = 1867.75
$ws.Cells.Item(664, 7).Value = 1112.599975585938
$ws.Cells.Item(664, 8).Value = 8507.349853515625
$ws.Cells.Item(664, 9).Value = -0.001637097066476197
$ws.Cells.Item(664, 10).Value = 174.2157176266467

$ws.Cells.Item(665, 1).Value = "'2024-09-12"
$ws.Cells.Item(665, 3).Value = 2247.5
$ws.Cells.Item(665, 4).Value = 1592.849975585938
$ws.Cells.Item(665, 5).Value = 1747.949951171875
$ws.Cells.Item(665, 6).Value = 1883.349975585938
$ws.Cells.Item(665, 7).Value = 1120.099975585938
$ws.Cells.Item(665, 8).Value = 8591.749877929688
$ws.Cells.Item(665, 9).Value = 0.009920836202496664
$ws.Cells.Item(665, 10).Value = 175.9440832251211

$ws.Cells.Item(666, 1).Value = "'2024-09-13"
$ws.Cells.Item(666, 3).Value = 2256.449951171875
$ws.Cells.Item(666, 4).Value = 1582.5
$ws.Cells.Item(666, 5).Value = 1753.699951171875
$ws.Cells.Item(666, 6).Value = 1923.300048828125
$ws.Cells.Item(666, 7).Value = 1118.550048828125
$ws.Cells.Item(666, 8).Value = 8634.5
$ws.Cells.Item(666, 9).Value = 0.004975717714982386
$ws.Cells.Item(666, 10).Value = 176.8195313168707

$ws.Cells.Item(667, 1).Value = "'2024-09-16"
$ws.Cells.Item(667, 3).Value = 2251.85009765625
$ws.Cells.Item(667, 4).Value = 1577.75
$ws.Cells.Item(667, 5).Value = 1741.449951171875
$ws.Cells.Item(667, 6).Value = 1900.949951171875
$ws.Cells.Item(667, 7).Value = 1115.849975585938
$ws.Cells.Item(667, 8).Value = 8587.849975585938
$ws.Cells.Item(667, 9).Value = -0.005402747630327465
$ws.Cells.Item(667, 10).Value = 175.8642200130528

$ws.Cells.Item(668, 1).Value = "'2024-09-17"
$ws.Cells.Item(668, 3).Value = 2270.39990234375
$ws.Cells.Item(668, 4).Value = 1561.699951171875
$ws.Cells.Item(668, 5).Value = 1713
$ws.Cells.Item(668, 6).Value = 1875.599975585938
$ws.Cells.Item(668, 7).Value = 1110.949951171875
$ws.Cells.Item(668, 8).Value = 8531.649780273438
$ws.Cells.Item(668, 9).Value = -0.006544151967287428
$ws.Cells.Item(668, 10).Value = 174.713337831679

$ws.Cells.Item(669, 1).Value = "'2024-09-18"
$ws.Cells.Item(669, 3).Value = 2224.949951171875
$ws.Cells.Item(669, 4).Value = 1543.050048828125
$ws.Cells.Item(669, 5).Value = 1646.050048828125
$ws.Cells.Item(669, 6).Value = 1857
$ws.Cells.Item(669, 7).Value = 1079.949951171875
$ws.Cells.Item(669, 8).Value = 8351
$ws.Cells.Item(669, 9).Value = -0.02117407358786916
$ws.Cells.Item(669, 10).Value = 171.0139447596488

$ws.Cells.Item(670, 1).Value = "'2024-09-19"
$ws.Cells.Item(670, 3).Value = 2171.89990234375
$ws.Cells.Item(670, 4).Value = 1515.050048828125
$ws.Cells.Item(670, 5).Value = 1649.800048828125
$ws.Cells.Item(670, 6).Value = 1886.5
$ws.Cells.Item(670, 7).Value = 1054.449951171875
$ws.Cells.Item(670, 8).Value = 8277.699951171875
$ws.Cells.Item(670, 9).Value = -0.00877739777608969
$ws.Cells.Item(670, 10).Value = 169.5128873412351

$ws.Cells.Item(671, 1).Value = "'2024-09-20"
$ws.Cells.Item(671, 3).Value = 2151.699951171875
$ws.Cells.Item(671, 4).Value = 1481.099975585938
$ws.Cells.Item(671, 5).Value = 1636.75
$ws.Cells.Item(671, 6).Value = 1897.25
$ws.Cells.Item(671, 7).Value = 1054.599975585938
$ws.Cells.Item(671, 8).Value = 8221.39990234375
$ws.Cells.Item(671, 9).Value = -0.0068014121265841
$ws.Cells.Item(671, 10).Value = 168.3599603336601

$ws.Cells.Item(672, 1).Value = "'2024-09-23"
$ws.Cells.Item(672, 3).Value = 2182.25
$ws.Cells.Item(672, 4).Value = 1440.400024414062
$ws.Cells.Item(672, 5).Value = 1712.449951171875
$ws.Cells.Item(672, 6).Value = 1952
$ws.Cells.Item(672, 7).Value = 1055.25
$ws.Cells.Item(672, 8).Value = 8342.349975585938
$ws.Cells.Item(672, 9).Value = 0.01471161537923817
$ws.Cells.Item(672, 10).Value = 170.8368073153527

$ws.Cells.Item(673, 1).Value = "'2024-09-24"
$ws.Cells.Item(673, 3).Value = 2215.75
$ws.Cells.Item(673, 4).Value = 1414.25
$ws.Cells.Item(673, 5).Value = 1697.5
$ws.Cells.Item(673, 6).Value = 1944.349975585938
$ws.Cells.Item(673, 7).Value = 1051.550048828125
$ws.Cells.Item(673, 8).Value = 8323.400024414062
$ws.Cells.Item(673, 9).Value = -0.002271536344954651
$ws.Cells.Item(673, 10).Value = 170.4487452984799

$ws.Cells.Item(674, 1).Value = "'2024-09-25"
$ws.Cells.Item(674, 3).Value = 2221.10009765625
$ws.Cells.Item(674, 4).Value = 1416.400024414062
$ws.Cells.Item(674, 5).Value = 1689.199951171875
$ws.Cells.Item(674, 6).Value = 1909.550048828125
$ws.Cells.Item(674, 7).Value = 1063.449951171875
$ws.Cells.Item(674, 8).Value = 8299.700073242188
$ws.Cells.Item(674, 9).Value = -0.002847388219040138
$ws.Cells.Item(674, 10).Value = 169.9634115491668

$ws.Cells.Item(675, 1).Value = "'2024-09-26"
$ws.Cells.Item(675, 3).Value = 2188.14990234375
$ws.Cells.Item(675, 4).Value = 1382.900024414062
$ws.Cells.Item(675, 5).Value = 1678.300048828125
$ws.Cells.Item(675, 6).Value = 1888.550048828125
$ws.Cells.Item(675, 7).Value = 1068
$ws.Cells.Item(675, 8).Value = 8205.900024414062
$ws.Cells.Item(675, 9).Value = -0.01130161909471061
$ws.Cells.Item(675, 10).Value = 168.0425498118006

$ws.Cells.Item(676, 1).Value = "'2024-09-27"
$ws.Cells.Item(676, 3).Value = 2218.5
$ws.Cells.Item(676, 4).Value = 1413.150024414062
$ws.Cells.Item(676, 5).Value = 1685.699951171875
$ws.Cells.Item(676, 6).Value = 1837.949951171875
$ws.Cells.Item(676, 7).Value = 1075.949951171875
$ws.Cells.Item(676, 8).Value = 8231.249877929688
$ws.Cells.Item(676, 9).Value = 0.003089222808004548
$ws.Cells.Item(676, 10).Value = 168.5616706893945
